$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1386.807
$ws.Range("J17").Value = 1411.3182
$ws.Range("L17").Value = 4233.9546
$ws.Range("N17").Value = -4569.9546
# Row 40
$ws.Range("H40").Value = 1285.8572
$ws.Range("I40").Value = 1067
$ws.Range("J40").Value = 1450
$ws.Range("K40").Value = 1067
$ws.Range("L40").Value = 1450
$ws.Range("M40").Value = -892
$ws.Range("N40").Value = -1800

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 977.55554
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 966.3333
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 966.3333
$ws.Range("M2").Value = -887
$ws.Range("N2").Value = -1192.3333
# Row 63
$ws.Range("H63").Value = 4619863.5
$ws.Range("I63").Value = 8149695
$ws.Range("J63").Value = 3930.7693
$ws.Range("K63").Value = 8149695
$ws.Range("L63").Value = 3930.7693
$ws.Range("M63").Value = -8149009
$ws.Range("N63").Value = -5302.7693
# Row 66
$ws.Range("H66").Value = 4619863.5
$ws.Range("I66").Value = 8149695
$ws.Range("J66").Value = 3930.7693
$ws.Range("K66").Value = 40748475
$ws.Range("L66").Value = 19653.8465
$ws.Range("M66").Value = -40745043
$ws.Range("N66").Value = -26517.8465
# Row 92
$ws.Range("H92").Value = 37999.5
$ws.Range("J92").Value = 37999.5
$ws.Range("L92").Value = 37999.5
$ws.Range("N92").Value = -42991.5
# Row 97
$ws.Range("H97").Value = 397.69232
$ws.Range("I97").Value = 428.09525
$ws.Range("J97").Value = 270
$ws.Range("K97").Value = 428.09525
$ws.Range("L97").Value = 270
$ws.Range("M97").Value = 67.90474999999998
$ws.Range("N97").Value = -1262
# Row 110
$ws.Range("H110").Value = 843.8214
$ws.Range("I110").Value = 825.5714
$ws.Range("J110").Value = 898.5714
$ws.Range("K110").Value = 825.5714
$ws.Range("L110").Value = 898.5714
$ws.Range("M110").Value = 1219.4286
$ws.Range("N110").Value = -4988.5714
# Row 116
$ws.Range("H116").Value = 977.55554
$ws.Range("I116").Value = 1000
$ws.Range("J116").Value = 966.3333
$ws.Range("K116").Value = 1000
$ws.Range("L116").Value = 966.3333
$ws.Range("M116").Value = 1294
$ws.Range("N116").Value = -5554.3333
# Row 122
$ws.Range("H122").Value = 1849.6666
$ws.Range("I122").Value = 1125.3334
$ws.Range("J122").Value = 3298.3333
$ws.Range("K122").Value = 3376.0002
$ws.Range("L122").Value = 9894.999899999999
$ws.Range("M122").Value = -926.0001999999999
$ws.Range("N122").Value = -14794.9999
# Row 132
$ws.Range("H132").Value = 2297.2144
$ws.Range("I132").Value = 1590.4615
$ws.Range("J132").Value = 3918.5881
$ws.Range("K132").Value = 4771.3845
$ws.Range("L132").Value = 11755.7643
$ws.Range("M132").Value = -2241.3845
$ws.Range("N132").Value = -16815.7643

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 977.55554
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 966.3333
$ws.Range("K3").Value = 1000
$ws.Range("L3").Value = 966.3333
$ws.Range("M3").Value = -886
$ws.Range("N3").Value = -1194.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1782.5256
$ws.Range("I58").Value = 1609.9688
$ws.Range("J58").Value = 2571.3572
$ws.Range("K58").Value = 1609.9688
$ws.Range("L58").Value = 2571.3572
$ws.Range("M58").Value = -1406.9688
$ws.Range("N58").Value = -2977.3572
# Row 62
$ws.Range("H62").Value = 100003740
$ws.Range("I62").Value = 100003740
$ws.Range("K62").Value = 100003740
$ws.Range("M62").Value = -100003116
# Row 65
$ws.Range("H65").Value = 100003740
$ws.Range("I65").Value = 100003740
$ws.Range("K65").Value = 500018700
$ws.Range("M65").Value = -500015580
# Row 86
$ws.Range("H86").Value = 3666.6667
$ws.Range("J86").Value = 3666.6667
$ws.Range("L86").Value = 3666.6667
$ws.Range("N86").Value = -5912.6667
# Row 89
$ws.Range("H89").Value = 3666.6667
$ws.Range("J89").Value = 3666.6667
$ws.Range("L89").Value = 18333.3335
$ws.Range("N89").Value = -29565.3335
# Row 99
$ws.Range("H99").Value = 20006272
$ws.Range("I99").Value = 50003950
$ws.Range("J99").Value = 7819
$ws.Range("K99").Value = 50003950
$ws.Range("L99").Value = 7819
$ws.Range("M99").Value = -50002452
$ws.Range("N99").Value = -10815
# Row 105
$ws.Range("H105").Value = 1556.95
$ws.Range("I105").Value = 1076.6
$ws.Range("J105").Value = 2998
$ws.Range("K105").Value = 1076.6
$ws.Range("L105").Value = 2998
$ws.Range("M105").Value = 670.4000000000001
$ws.Range("N105").Value = -6492
# Row 126
$ws.Range("H126").Value = 20006272
$ws.Range("I126").Value = 50003950
$ws.Range("J126").Value = 7819
$ws.Range("K126").Value = 150011850
$ws.Range("L126").Value = 23457
$ws.Range("M126").Value = -150009380
$ws.Range("N126").Value = -28397
# Row 134
$ws.Range("H134").Value = 5813
$ws.Range("I134").Value = 6581.778
$ws.Range("K134").Value = 19745.334
$ws.Range("M134").Value = -17210.334
# Row 136
$ws.Range("H136").Value = 1782.5256
$ws.Range("I136").Value = 1609.9688
$ws.Range("J136").Value = 2571.3572
$ws.Range("K136").Value = 4829.9064
$ws.Range("L136").Value = 7714.071599999999
$ws.Range("M136").Value = -2279.9064
$ws.Range("N136").Value = -12814.0716

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 93
$ws.Range("H93").Value = 6006.75
$ws.Range("J93").Value = 6006.75
$ws.Range("L93").Value = 18020.25
$ws.Range("N93").Value = -21764.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 11
$ws.Range("H11").Value = 8457300
$ws.Range("I11").Value = 17500000
$ws.Range("J11").Value = 3290043.5
$ws.Range("K11").Value = 17500000
$ws.Range("L11").Value = 3290043.5
$ws.Range("M11").Value = -17499861
$ws.Range("N11").Value = -3290321.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 9649.532999999999
$ws.Range("I40").Value = 8995.111000000001
$ws.Range("J40").Value = 10631.167
$ws.Range("K40").Value = 8995.111000000001
$ws.Range("L40").Value = 10631.167
$ws.Range("M40").Value = -8859.111000000001
$ws.Range("N40").Value = -10903.167
# Row 61
$ws.Range("H61").Value = 1435.8636
$ws.Range("I61").Value = 1417
$ws.Range("K61").Value = 1417
$ws.Range("M61").Value = -1215
# Row 113
$ws.Range("H113").Value = 1435.8636
$ws.Range("I113").Value = 1417
$ws.Range("K113").Value = 1417
$ws.Range("M113").Value = 753

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 111168220
$ws.Range("I96").Value = 333500670
$ws.Range("J96").Value = 1999
$ws.Range("K96").Value = 333500670
$ws.Range("L96").Value = 1999
$ws.Range("M96").Value = -333499297
$ws.Range("N96").Value = -4745
# Row 113
$ws.Range("H113").Value = 449.41666
$ws.Range("I113").Value = 377
$ws.Range("K113").Value = 1131
$ws.Range("M113").Value = 1039
